$d = $word.ActiveDocument

# Phase 1: replace each original expression with a unique placeholder
# to avoid collisions where a new value is a substring of a later old value.
$null = $d.Content.Find.Execute("74-37=", $true, $false, $false, $false, $false, $true, 1, $false, "@@0@@", 2)
$null = $d.Content.Find.Execute("31-9=", $true, $false, $false, $false, $false, $true, 1, $false, "@@1@@", 2)
$null = $d.Content.Find.Execute("0+20=", $true, $false, $false, $false, $false, $true, 1, $false, "@@2@@", 2)
$null = $d.Content.Find.Execute("67-47=", $true, $false, $false, $false, $false, $true, 1, $false, "@@3@@", 2)
$null = $d.Content.Find.Execute("77-11=", $true, $false, $false, $false, $false, $true, 1, $false, "@@4@@", 2)
$null = $d.Content.Find.Execute("89-81=", $true, $false, $false, $false, $false, $true, 1, $false, "@@5@@", 2)
$null = $d.Content.Find.Execute("9+57=", $true, $false, $false, $false, $false, $true, 1, $false, "@@6@@", 2)
$null = $d.Content.Find.Execute("21+57=", $true, $false, $false, $false, $false, $true, 1, $false, "@@7@@", 2)
$null = $d.Content.Find.Execute("22+67=", $true, $false, $false, $false, $false, $true, 1, $false, "@@8@@", 2)
$null = $d.Content.Find.Execute("90+4=", $true, $false, $false, $false, $false, $true, 1, $false, "@@9@@", 2)
$null = $d.Content.Find.Execute("86-10=", $true, $false, $false, $false, $false, $true, 1, $false, "@@10@@", 2)
$null = $d.Content.Find.Execute("68-1=", $true, $false, $false, $false, $false, $true, 1, $false, "@@11@@", 2)
$null = $d.Content.Find.Execute("15+7=", $true, $false, $false, $false, $false, $true, 1, $false, "@@12@@", 2)
$null = $d.Content.Find.Execute("97-4=", $true, $false, $false, $false, $false, $true, 1, $false, "@@13@@", 2)
$null = $d.Content.Find.Execute("86-0=", $true, $false, $false, $false, $false, $true, 1, $false, "@@14@@", 2)
$null = $d.Content.Find.Execute("70+6=", $true, $false, $false, $false, $false, $true, 1, $false, "@@15@@", 2)
$null = $d.Content.Find.Execute("44+15=", $true, $false, $false, $false, $false, $true, 1, $false, "@@16@@", 2)
$null = $d.Content.Find.Execute("41-28=", $true, $false, $false, $false, $false, $true, 1, $false, "@@17@@", 2)
$null = $d.Content.Find.Execute("29+51=", $true, $false, $false, $false, $false, $true, 1, $false, "@@18@@", 2)
$null = $d.Content.Find.Execute("75-53=", $true, $false, $false, $false, $false, $true, 1, $false, "@@19@@", 2)
$null = $d.Content.Find.Execute("0+61=", $true, $false, $false, $false, $false, $true, 1, $false, "@@20@@", 2)
$null = $d.Content.Find.Execute("62-26=", $true, $false, $false, $false, $false, $true, 1, $false, "@@21@@", 2)
$null = $d.Content.Find.Execute("63-38=", $true, $false, $false, $false, $false, $true, 1, $false, "@@22@@", 2)
$null = $d.Content.Find.Execute("80+2=", $true, $false, $false, $false, $false, $true, 1, $false, "@@23@@", 2)
$null = $d.Content.Find.Execute("9-7=", $true, $false, $false, $false, $false, $true, 1, $false, "@@24@@", 2)
$null = $d.Content.Find.Execute("80+0=", $true, $false, $false, $false, $false, $true, 1, $false, "@@25@@", 2)
$null = $d.Content.Find.Execute("9+60=", $true, $false, $false, $false, $false, $true, 1, $false, "@@26@@", 2)
$null = $d.Content.Find.Execute("77-12=", $true, $false, $false, $false, $false, $true, 1, $false, "@@27@@", 2)
$null = $d.Content.Find.Execute("86-67=", $true, $false, $false, $false, $false, $true, 1, $false, "@@28@@", 2)
$null = $d.Content.Find.Execute("61-52=", $true, $false, $false, $false, $false, $true, 1, $false, "@@29@@", 2)
$null = $d.Content.Find.Execute("78+9=", $true, $false, $false, $false, $false, $true, 1, $false, "@@30@@", 2)
$null = $d.Content.Find.Execute("91-41=", $true, $false, $false, $false, $false, $true, 1, $false, "@@31@@", 2)
$null = $d.Content.Find.Execute("55-35=", $true, $false, $false, $false, $false, $true, 1, $false, "@@32@@", 2)
$null = $d.Content.Find.Execute("26+48=", $true, $false, $false, $false, $false, $true, 1, $false, "@@33@@", 2)
$null = $d.Content.Find.Execute("20+8=", $true, $false, $false, $false, $false, $true, 1, $false, "@@34@@", 2)
$null = $d.Content.Find.Execute("71-31=", $true, $false, $false, $false, $false, $true, 1, $false, "@@35@@", 2)
$null = $d.Content.Find.Execute("8+63=", $true, $false, $false, $false, $false, $true, 1, $false, "@@36@@", 2)
$null = $d.Content.Find.Execute("51+3=", $true, $false, $false, $false, $false, $true, 1, $false, "@@37@@", 2)
$null = $d.Content.Find.Execute("45-20=", $true, $false, $false, $false, $false, $true, 1, $false, "@@38@@", 2)
$null = $d.Content.Find.Execute("38+32=", $true, $false, $false, $false, $false, $true, 1, $false, "@@39@@", 2)
$null = $d.Content.Find.Execute("7+56=", $true, $false, $false, $false, $false, $true, 1, $false, "@@40@@", 2)
$null = $d.Content.Find.Execute("16+40=", $true, $false, $false, $false, $false, $true, 1, $false, "@@41@@", 2)
$null = $d.Content.Find.Execute("4+74=", $true, $false, $false, $false, $false, $true, 1, $false, "@@42@@", 2)
$null = $d.Content.Find.Execute("6+76=", $true, $false, $false, $false, $false, $true, 1, $false, "@@43@@", 2)
$null = $d.Content.Find.Execute("38-1=", $true, $false, $false, $false, $false, $true, 1, $false, "@@44@@", 2)
$null = $d.Content.Find.Execute("8+13=", $true, $false, $false, $false, $false, $true, 1, $false, "@@45@@", 2)
$null = $d.Content.Find.Execute("28-4=", $true, $false, $false, $false, $false, $true, 1, $false, "@@46@@", 2)
$null = $d.Content.Find.Execute("75-16=", $true, $false, $false, $false, $false, $true, 1, $false, "@@47@@", 2)
$null = $d.Content.Find.Execute("34+1=", $true, $false, $false, $false, $false, $true, 1, $false, "@@48@@", 2)
$null = $d.Content.Find.Execute("37-26=", $true, $false, $false, $false, $false, $true, 1, $false, "@@49@@", 2)
$null = $d.Content.Find.Execute("16+32=", $true, $false, $false, $false, $false, $true, 1, $false, "@@50@@", 2)
$null = $d.Content.Find.Execute("39+38=", $true, $false, $false, $false, $false, $true, 1, $false, "@@51@@", 2)
$null = $d.Content.Find.Execute("41+54=", $true, $false, $false, $false, $false, $true, 1, $false, "@@52@@", 2)
$null = $d.Content.Find.Execute("74-7=", $true, $false, $false, $false, $false, $true, 1, $false, "@@53@@", 2)
$null = $d.Content.Find.Execute("98-76=", $true, $false, $false, $false, $false, $true, 1, $false, "@@54@@", 2)
$null = $d.Content.Find.Execute("10+75=", $true, $false, $false, $false, $false, $true, 1, $false, "@@55@@", 2)
$null = $d.Content.Find.Execute("59-12=", $true, $false, $false, $false, $false, $true, 1, $false, "@@56@@", 2)
$null = $d.Content.Find.Execute("68-53=", $true, $false, $false, $false, $false, $true, 1, $false, "@@57@@", 2)
$null = $d.Content.Find.Execute("99-86=", $true, $false, $false, $false, $false, $true, 1, $false, "@@58@@", 2)
$null = $d.Content.Find.Execute("51-25=", $true, $false, $false, $false, $false, $true, 1, $false, "@@59@@", 2)
$null = $d.Content.Find.Execute("20+17=", $true, $false, $false, $false, $false, $true, 1, $false, "@@60@@", 2)
$null = $d.Content.Find.Execute("42+8=", $true, $false, $false, $false, $false, $true, 1, $false, "@@61@@", 2)
$null = $d.Content.Find.Execute("93-77=", $true, $false, $false, $false, $false, $true, 1, $false, "@@62@@", 2)
$null = $d.Content.Find.Execute("40+22=", $true, $false, $false, $false, $false, $true, 1, $false, "@@63@@", 2)
$null = $d.Content.Find.Execute("97-73=", $true, $false, $false, $false, $false, $true, 1, $false, "@@64@@", 2)
$null = $d.Content.Find.Execute("81-10=", $true, $false, $false, $false, $false, $true, 1, $false, "@@65@@", 2)
$null = $d.Content.Find.Execute("63+9=", $true, $false, $false, $false, $false, $true, 1, $false, "@@66@@", 2)
$null = $d.Content.Find.Execute("39+14=", $true, $false, $false, $false, $false, $true, 1, $false, "@@67@@", 2)
$null = $d.Content.Find.Execute("0+87=", $true, $false, $false, $false, $false, $true, 1, $false, "@@68@@", 2)
$null = $d.Content.Find.Execute("96-76=", $true, $false, $false, $false, $false, $true, 1, $false, "@@69@@", 2)
$null = $d.Content.Find.Execute("99-89=", $true, $false, $false, $false, $false, $true, 1, $false, "@@70@@", 2)
$null = $d.Content.Find.Execute("9+8=", $true, $false, $false, $false, $false, $true, 1, $false, "@@71@@", 2)
$null = $d.Content.Find.Execute("85-74=", $true, $false, $false, $false, $false, $true, 1, $false, "@@72@@", 2)
$null = $d.Content.Find.Execute("36-1=", $true, $false, $false, $false, $false, $true, 1, $false, "@@73@@", 2)
$null = $d.Content.Find.Execute("91-88=", $true, $false, $false, $false, $false, $true, 1, $false, "@@74@@", 2)
$null = $d.Content.Find.Execute("2+8=", $true, $false, $false, $false, $false, $true, 1, $false, "@@75@@", 2)
$null = $d.Content.Find.Execute("30+37=", $true, $false, $false, $false, $false, $true, 1, $false, "@@76@@", 2)
$null = $d.Content.Find.Execute("76-75=", $true, $false, $false, $false, $false, $true, 1, $false, "@@77@@", 2)
$null = $d.Content.Find.Execute("2+54=", $true, $false, $false, $false, $false, $true, 1, $false, "@@78@@", 2)
$null = $d.Content.Find.Execute("9+39=", $true, $false, $false, $false, $false, $true, 1, $false, "@@79@@", 2)
$null = $d.Content.Find.Execute("94-29=", $true, $false, $false, $false, $false, $true, 1, $false, "@@80@@", 2)
$null = $d.Content.Find.Execute("81+9=", $true, $false, $false, $false, $false, $true, 1, $false, "@@81@@", 2)
$null = $d.Content.Find.Execute("24+62=", $true, $false, $false, $false, $false, $true, 1, $false, "@@82@@", 2)
$null = $d.Content.Find.Execute("70-44=", $true, $false, $false, $false, $false, $true, 1, $false, "@@83@@", 2)
$null = $d.Content.Find.Execute("6-5=", $true, $false, $false, $false, $false, $true, 1, $false, "@@84@@", 2)
$null = $d.Content.Find.Execute("78-60=", $true, $false, $false, $false, $false, $true, 1, $false, "@@85@@", 2)
$null = $d.Content.Find.Execute("33+24=", $true, $false, $false, $false, $false, $true, 1, $false, "@@86@@", 2)
$null = $d.Content.Find.Execute("20+51=", $true, $false, $false, $false, $false, $true, 1, $false, "@@87@@", 2)
$null = $d.Content.Find.Execute("85-64=", $true, $false, $false, $false, $false, $true, 1, $false, "@@88@@", 2)
$null = $d.Content.Find.Execute("35+38=", $true, $false, $false, $false, $false, $true, 1, $false, "@@89@@", 2)
$null = $d.Content.Find.Execute("54+13=", $true, $false, $false, $false, $false, $true, 1, $false, "@@90@@", 2)
$null = $d.Content.Find.Execute("81-49=", $true, $false, $false, $false, $false, $true, 1, $false, "@@91@@", 2)
$null = $d.Content.Find.Execute("73-1=", $true, $false, $false, $false, $false, $true, 1, $false, "@@92@@", 2)
$null = $d.Content.Find.Execute("35+8=", $true, $false, $false, $false, $false, $true, 1, $false, "@@93@@", 2)
$null = $d.Content.Find.Execute("36+26=", $true, $false, $false, $false, $false, $true, 1, $false, "@@94@@", 2)
$null = $d.Content.Find.Execute("71+27=", $true, $false, $false, $false, $false, $true, 1, $false, "@@95@@", 2)
$null = $d.Content.Find.Execute("71+26=", $true, $false, $false, $false, $false, $true, 1, $false, "@@96@@", 2)
$null = $d.Content.Find.Execute("77+22=", $true, $false, $false, $false, $false, $true, 1, $false, "@@97@@", 2)
$null = $d.Content.Find.Execute("71-4=", $true, $false, $false, $false, $false, $true, 1, $false, "@@98@@", 2)
$null = $d.Content.Find.Execute("38-0=", $true, $false, $false, $false, $false, $true, 1, $false, "@@99@@", 2)

# Phase 2: replace each placeholder with the final new expression
$null = $d.Content.Find.Execute("@@0@@", $true, $false, $false, $false, $false, $true, 1, $false, "7+86=", 2)
$null = $d.Content.Find.Execute("@@1@@", $true, $false, $false, $false, $false, $true, 1, $false, "86-56=", 2)
$null = $d.Content.Find.Execute("@@2@@", $true, $false, $false, $false, $false, $true, 1, $false, "8-7=", 2)
$null = $d.Content.Find.Execute("@@3@@", $true, $false, $false, $false, $false, $true, 1, $false, "66+5=", 2)
$null = $d.Content.Find.Execute("@@4@@", $true, $false, $false, $false, $false, $true, 1, $false, "39-25=", 2)
$null = $d.Content.Find.Execute("@@5@@", $true, $false, $false, $false, $false, $true, 1, $false, "78-8=", 2)
$null = $d.Content.Find.Execute("@@6@@", $true, $false, $false, $false, $false, $true, 1, $false, "13+52=", 2)
$null = $d.Content.Find.Execute("@@7@@", $true, $false, $false, $false, $false, $true, 1, $false, "19+61=", 2)
$null = $d.Content.Find.Execute("@@8@@", $true, $false, $false, $false, $false, $true, 1, $false, "48+0=", 2)
$null = $d.Content.Find.Execute("@@9@@", $true, $false, $false, $false, $false, $true, 1, $false, "57-50=", 2)
$null = $d.Content.Find.Execute("@@10@@", $true, $false, $false, $false, $false, $true, 1, $false, "33-25=", 2)
$null = $d.Content.Find.Execute("@@11@@", $true, $false, $false, $false, $false, $true, 1, $false, "48-47=", 2)
$null = $d.Content.Find.Execute("@@12@@", $true, $false, $false, $false, $false, $true, 1, $false, "24+29=", 2)
$null = $d.Content.Find.Execute("@@13@@", $true, $false, $false, $false, $false, $true, 1, $false, "4+56=", 2)
$null = $d.Content.Find.Execute("@@14@@", $true, $false, $false, $false, $false, $true, 1, $false, "56+4=", 2)
$null = $d.Content.Find.Execute("@@15@@", $true, $false, $false, $false, $false, $true, 1, $false, "59-4=", 2)
$null = $d.Content.Find.Execute("@@16@@", $true, $false, $false, $false, $false, $true, 1, $false, "50+18=", 2)
$null = $d.Content.Find.Execute("@@17@@", $true, $false, $false, $false, $false, $true, 1, $false, "16+46=", 2)
$null = $d.Content.Find.Execute("@@18@@", $true, $false, $false, $false, $false, $true, 1, $false, "95-22=", 2)
$null = $d.Content.Find.Execute("@@19@@", $true, $false, $false, $false, $false, $true, 1, $false, "73+25=", 2)
$null = $d.Content.Find.Execute("@@20@@", $true, $false, $false, $false, $false, $true, 1, $false, "30+10=", 2)
$null = $d.Content.Find.Execute("@@21@@", $true, $false, $false, $false, $false, $true, 1, $false, "19+1=", 2)
$null = $d.Content.Find.Execute("@@22@@", $true, $false, $false, $false, $false, $true, 1, $false, "86-43=", 2)
$null = $d.Content.Find.Execute("@@23@@", $true, $false, $false, $false, $false, $true, 1, $false, "42-26=", 2)
$null = $d.Content.Find.Execute("@@24@@", $true, $false, $false, $false, $false, $true, 1, $false, "83-59=", 2)
$null = $d.Content.Find.Execute("@@25@@", $true, $false, $false, $false, $false, $true, 1, $false, "61-2=", 2)
$null = $d.Content.Find.Execute("@@26@@", $true, $false, $false, $false, $false, $true, 1, $false, "35+36=", 2)
$null = $d.Content.Find.Execute("@@27@@", $true, $false, $false, $false, $false, $true, 1, $false, "60-34=", 2)
$null = $d.Content.Find.Execute("@@28@@", $true, $false, $false, $false, $false, $true, 1, $false, "41+33=", 2)
$null = $d.Content.Find.Execute("@@29@@", $true, $false, $false, $false, $false, $true, 1, $false, "84-51=", 2)
$null = $d.Content.Find.Execute("@@30@@", $true, $false, $false, $false, $false, $true, 1, $false, "32+54=", 2)
$null = $d.Content.Find.Execute("@@31@@", $true, $false, $false, $false, $false, $true, 1, $false, "40-17=", 2)
$null = $d.Content.Find.Execute("@@32@@", $true, $false, $false, $false, $false, $true, 1, $false, "57+24=", 2)
$null = $d.Content.Find.Execute("@@33@@", $true, $false, $false, $false, $false, $true, 1, $false, "75+0=", 2)
$null = $d.Content.Find.Execute("@@34@@", $true, $false, $false, $false, $false, $true, 1, $false, "58-54=", 2)
$null = $d.Content.Find.Execute("@@35@@", $true, $false, $false, $false, $false, $true, 1, $false, "81-23=", 2)
$null = $d.Content.Find.Execute("@@36@@", $true, $false, $false, $false, $false, $true, 1, $false, "41+0=", 2)
$null = $d.Content.Find.Execute("@@37@@", $true, $false, $false, $false, $false, $true, 1, $false, "80-49=", 2)
$null = $d.Content.Find.Execute("@@38@@", $true, $false, $false, $false, $false, $true, 1, $false, "65-58=", 2)
$null = $d.Content.Find.Execute("@@39@@", $true, $false, $false, $false, $false, $true, 1, $false, "3+79=", 2)
$null = $d.Content.Find.Execute("@@40@@", $true, $false, $false, $false, $false, $true, 1, $false, "11+47=", 2)
$null = $d.Content.Find.Execute("@@41@@", $true, $false, $false, $false, $false, $true, 1, $false, "27-20=", 2)
$null = $d.Content.Find.Execute("@@42@@", $true, $false, $false, $false, $false, $true, 1, $false, "22+55=", 2)
$null = $d.Content.Find.Execute("@@43@@", $true, $false, $false, $false, $false, $true, 1, $false, "82-72=", 2)
$null = $d.Content.Find.Execute("@@44@@", $true, $false, $false, $false, $false, $true, 1, $false, "82-47=", 2)
$null = $d.Content.Find.Execute("@@45@@", $true, $false, $false, $false, $false, $true, 1, $false, "38-34=", 2)
$null = $d.Content.Find.Execute("@@46@@", $true, $false, $false, $false, $false, $true, 1, $false, "7+23=", 2)
$null = $d.Content.Find.Execute("@@47@@", $true, $false, $false, $false, $false, $true, 1, $false, "92-71=", 2)
$null = $d.Content.Find.Execute("@@48@@", $true, $false, $false, $false, $false, $true, 1, $false, "18+50=", 2)
$null = $d.Content.Find.Execute("@@49@@", $true, $false, $false, $false, $false, $true, 1, $false, "97-21=", 2)
$null = $d.Content.Find.Execute("@@50@@", $true, $false, $false, $false, $false, $true, 1, $false, "83+4=", 2)
$null = $d.Content.Find.Execute("@@51@@", $true, $false, $false, $false, $false, $true, 1, $false, "54-28=", 2)
$null = $d.Content.Find.Execute("@@52@@", $true, $false, $false, $false, $false, $true, 1, $false, "4+86=", 2)
$null = $d.Content.Find.Execute("@@53@@", $true, $false, $false, $false, $false, $true, 1, $false, "37-32=", 2)
$null = $d.Content.Find.Execute("@@54@@", $true, $false, $false, $false, $false, $true, 1, $false, "64+19=", 2)
$null = $d.Content.Find.Execute("@@55@@", $true, $false, $false, $false, $false, $true, 1, $false, "20+40=", 2)
$null = $d.Content.Find.Execute("@@56@@", $true, $false, $false, $false, $false, $true, 1, $false, "63-30=", 2)
$null = $d.Content.Find.Execute("@@57@@", $true, $false, $false, $false, $false, $true, 1, $false, "1+16=", 2)
$null = $d.Content.Find.Execute("@@58@@", $true, $false, $false, $false, $false, $true, 1, $false, "19-10=", 2)
$null = $d.Content.Find.Execute("@@59@@", $true, $false, $false, $false, $false, $true, 1, $false, "85+3=", 2)
$null = $d.Content.Find.Execute("@@60@@", $true, $false, $false, $false, $false, $true, 1, $false, "64+23=", 2)
$null = $d.Content.Find.Execute("@@61@@", $true, $false, $false, $false, $false, $true, 1, $false, "34+25=", 2)
$null = $d.Content.Find.Execute("@@62@@", $true, $false, $false, $false, $false, $true, 1, $false, "47+35=", 2)
$null = $d.Content.Find.Execute("@@63@@", $true, $false, $false, $false, $false, $true, 1, $false, "50-30=", 2)
$null = $d.Content.Find.Execute("@@64@@", $true, $false, $false, $false, $false, $true, 1, $false, "20-8=", 2)
$null = $d.Content.Find.Execute("@@65@@", $true, $false, $false, $false, $false, $true, 1, $false, "84-73=", 2)
$null = $d.Content.Find.Execute("@@66@@", $true, $false, $false, $false, $false, $true, 1, $false, "32+6=", 2)
$null = $d.Content.Find.Execute("@@67@@", $true, $false, $false, $false, $false, $true, 1, $false, "87+8=", 2)
$null = $d.Content.Find.Execute("@@68@@", $true, $false, $false, $false, $false, $true, 1, $false, "38+5=", 2)
$null = $d.Content.Find.Execute("@@69@@", $true, $false, $false, $false, $false, $true, 1, $false, "92-19=", 2)
$null = $d.Content.Find.Execute("@@70@@", $true, $false, $false, $false, $false, $true, 1, $false, "68-4=", 2)
$null = $d.Content.Find.Execute("@@71@@", $true, $false, $false, $false, $false, $true, 1, $false, "26-12=", 2)
$null = $d.Content.Find.Execute("@@72@@", $true, $false, $false, $false, $false, $true, 1, $false, "45+1=", 2)
$null = $d.Content.Find.Execute("@@73@@", $true, $false, $false, $false, $false, $true, 1, $false, "60+21=", 2)
$null = $d.Content.Find.Execute("@@74@@", $true, $false, $false, $false, $false, $true, 1, $false, "74-47=", 2)
$null = $d.Content.Find.Execute("@@75@@", $true, $false, $false, $false, $false, $true, 1, $false, "92-51=", 2)
$null = $d.Content.Find.Execute("@@76@@", $true, $false, $false, $false, $false, $true, 1, $false, "97-67=", 2)
$null = $d.Content.Find.Execute("@@77@@", $true, $false, $false, $false, $false, $true, 1, $false, "77-73=", 2)
$null = $d.Content.Find.Execute("@@78@@", $true, $false, $false, $false, $false, $true, 1, $false, "48+8=", 2)
$null = $d.Content.Find.Execute("@@79@@", $true, $false, $false, $false, $false, $true, 1, $false, "23-19=", 2)
$null = $d.Content.Find.Execute("@@80@@", $true, $false, $false, $false, $false, $true, 1, $false, "96-86=", 2)
$null = $d.Content.Find.Execute("@@81@@", $true, $false, $false, $false, $false, $true, 1, $false, "14+17=", 2)
$null = $d.Content.Find.Execute("@@82@@", $true, $false, $false, $false, $false, $true, 1, $false, "67+30=", 2)
$null = $d.Content.Find.Execute("@@83@@", $true, $false, $false, $false, $false, $true, 1, $false, "87-74=", 2)
$null = $d.Content.Find.Execute("@@84@@", $true, $false, $false, $false, $false, $true, 1, $false, "39-9=", 2)
$null = $d.Content.Find.Execute("@@85@@", $true, $false, $false, $false, $false, $true, 1, $false, "69-31=", 2)
$null = $d.Content.Find.Execute("@@86@@", $true, $false, $false, $false, $false, $true, 1, $false, "98-66=", 2)
$null = $d.Content.Find.Execute("@@87@@", $true, $false, $false, $false, $false, $true, 1, $false, "25+10=", 2)
$null = $d.Content.Find.Execute("@@88@@", $true, $false, $false, $false, $false, $true, 1, $false, "52+42=", 2)
$null = $d.Content.Find.Execute("@@89@@", $true, $false, $false, $false, $false, $true, 1, $false, "94-90=", 2)
$null = $d.Content.Find.Execute("@@90@@", $true, $false, $false, $false, $false, $true, 1, $false, "39+36=", 2)
$null = $d.Content.Find.Execute("@@91@@", $true, $false, $false, $false, $false, $true, 1, $false, "53+14=", 2)
$null = $d.Content.Find.Execute("@@92@@", $true, $false, $false, $false, $false, $true, 1, $false, "57-3=", 2)
$null = $d.Content.Find.Execute("@@93@@", $true, $false, $false, $false, $false, $true, 1, $false, "36-29=", 2)
$null = $d.Content.Find.Execute("@@94@@", $true, $false, $false, $false, $false, $true, 1, $false, "73+7=", 2)
$null = $d.Content.Find.Execute("@@95@@", $true, $false, $false, $false, $false, $true, 1, $false, "31-2=", 2)
$null = $d.Content.Find.Execute("@@96@@", $true, $false, $false, $false, $false, $true, 1, $false, "81-12=", 2)
$null = $d.Content.Find.Execute("@@97@@", $true, $false, $false, $false, $false, $true, 1, $false, "24+52=", 2)
$null = $d.Content.Find.Execute("@@98@@", $true, $false, $false, $false, $false, $true, 1, $false, "0+1=", 2)
$null = $d.Content.Find.Execute("@@99@@", $true, $false, $false, $false, $false, $true, 1, $false, "69-56=", 2)
